# Add a new "Stamp Image" table column between "Display Name" and "Status"
# on the User sheet, populate the two known stamp paths, and update the
# selected cell to match the post-edit workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$tbl = $ws.ListObjects.Item(1)

# Shift column E (and everything right of it) one column to the right so a
# new, empty column E is created; xlShiftToRight = -4161.
$null = $ws.Range("E1:E11").Insert(-4161)

# Grow the table to include the freshly inserted column.
$tbl.Resize($ws.Range("A1:F11"))

# Header row: new "Stamp Image" header in E1, "Status" header shifted to F1.
$ws.Cells.Item(1, 5).Value = "Stamp Image"
$ws.Cells.Item(1, 6).Value = "Status"

# Populate the two known user stamp image paths.
$ws.Cells.Item(2, 5).Value = "\STAMP\adib.jamil.PNG"
$ws.Cells.Item(7, 5).Value = "\STAMP\hakim.hisham.PNG"

# Match the new column's width.
$ws.Columns.Item(5).ColumnWidth = 36.7

# Update the active selection to D7.
$null = $ws.Range("D7").Select()
